# Applies the "Update works images 2026-01-18 07:01:36" edit:
#  1) Sheet1 ("names" list): the top id ("8a24swce") has been consumed / used,
#     so its row is removed and every following row shifts up by one.
#  2) "used" sheet: a new usage record is appended for that id, pointing at the
#     newly generated image file and its usage timestamp.

$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# Capture the id that is being "used" before it disappears from Sheet1.
$usedId = $namesSheet.Cells.Item(1, 1).Value()

# 1) Remove the consumed id from the top of the names list; rows below shift up.
$namesSheet.Rows.Item(1).Delete()

# 2) Append the new usage record to the "used" sheet.
$newRow = $usedSheet.Cells.Item($usedSheet.Rows.Count, 1).End(-4162).Row + 1

$usedSheet.Cells.Item($newRow, 1).Value = $usedId
$usedSheet.Cells.Item($newRow, 2).Value = "ChatGPT Image 2026年1月18日 07_00_29.png"
$usedSheet.Cells.Item($newRow, 3).Value = "2026-01-18 07:01:32"
